# Applies the "feat: using new template" edit to the PETALBURG CITY
# coordinates sheet:
#   - extends the tracked grid from column AL out to column AU (9 more
#     columns) on the header row and on a handful of data rows
#   - re-labels a batch of existing "1" marker cells on rows 5-10 and
#     15-20 as the text "farming" instead of the numeric 1
#   - moves the active selection to AW8
#
# Excel xlPasteFormats constant (Paste Special -> Formats only).
$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Header row: new sequence numbers 44..52 in AM1:AU1, styled like the
#    existing header cells (bold + border), copied from AL1.
# ---------------------------------------------------------------------
$headerCols = @("AM","AN","AO","AP","AQ","AR","AS","AT","AU")
$headerVals = @(44,45,46,47,48,49,50,51,52)
for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $ws.Range($headerCols[$i] + "1").Value2 = $headerVals[$i]
}
$ws.Range("AL1").Copy() | Out-Null
$ws.Range("AM1:AU1").PasteSpecial($xlPasteFormats) | Out-Null

# ---------------------------------------------------------------------
# 2) New "1" marker cells on rows 5-8, extending the existing pattern
#    out to column AU (style copied from the row's existing AL cell,
#    which already carries the data-cell style). Gaps (AP5, AR7) are
#    intentionally left untouched so no cell gets created there.
# ---------------------------------------------------------------------
function Fill-NewMarkers($row, $ranges) {
    $src = "AL" + $row
    foreach ($r in $ranges) {
        $rangeRef = $r[0] + $row + ":" + $r[1] + $row
        $ws.Range($rangeRef).Value2 = 1
        $ws.Range($src).Copy() | Out-Null
        $ws.Range($rangeRef).PasteSpecial($xlPasteFormats) | Out-Null
    }
}

Fill-NewMarkers 5 @(@("AM","AO"), @("AQ","AU"))
Fill-NewMarkers 6 @(@("AM","AU"))
Fill-NewMarkers 7 @(@("AM","AQ"), @("AS","AU"))
Fill-NewMarkers 8 @(@("AM","AU"))

# ---------------------------------------------------------------------
# 3) Re-label existing numeric "1" cells as the text "farming" (these
#    become shared-string cells; style/formatting is untouched).
# ---------------------------------------------------------------------
$farmingRanges = @(
    @(5,  "AG", "AI"),
    @(6,  "AG", "AI"),
    @(7,  "AA", "AI"),
    @(8,  "AA", "AH"),
    @(9,  "AB", "AH"),
    @(10, "AC", "AF"),
    @(15, "AE", "AF"),
    @(16, "AE", "AF"),
    @(17, "AE", "AG"),
    @(18, "AE", "AG"),
    @(19, "AE", "AH"),
    @(20, "AE", "AH")
)

foreach ($entry in $farmingRanges) {
    $row = $entry[0]
    $fromCol = $entry[1]
    $toCol = $entry[2]
    $rangeRef = $fromCol + $row + ":" + $toCol + $row
    $ws.Range($rangeRef).Value2 = "farming"
}

# ---------------------------------------------------------------------
# 4) Move the active selection to AW8 (matches the saved view state).
# ---------------------------------------------------------------------
$ws.Range("AW8").Select() | Out-Null
